$wb = $excel.ActiveWorkbook

# --- Fix the "2050" column header label (it previously held a stray
#     leftover numeric value instead of the intended year label) ---
# Sheets that use plain year labels ("2015"/"2030"/"2040"/"2050")
$yearSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)
foreach ($name in $yearSheets) {
    $ws = $wb.Worksheets.Item($name)
    # Force the cell to hold literal text "2050" (not the number 2050),
    # matching its neighbouring header cells (B1:D1), while keeping the
    # exact same cell formatting (bold/centered/bordered header style).
    $ws.Range("E1").NumberFormat = "@"
    $ws.Range("E1").Value = "2050"
    $ws.Range("D1").Copy()
    $ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
}

# Sheet 4 uses year-range labels instead of plain years
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("E1").NumberFormat = "@"
$ws4.Range("E1").Value = "2041-2050"
$ws4.Range("D1").Copy()
$ws4.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# --- Remove the "Total" summary rows from the bottom of each table ---
$totalRowSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)
foreach ($name in $totalRowSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(13).Delete()
}

# Sheet 6 has its Total row at row 4 instead of 13
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
